# Add materials for session 02 (row 3 of the schedule table)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Slides link for session 02
$ws.Range("E3").Value = "slides/slides.html#/sitzung-02-warum-wir-mediennutzung-unterhaltsam-finden"

# Exercise link for session 02
$ws.Range("F3").Value = "exercises/e02.html"

# Preparation reading link for session 03 (row 4)
$ws.Range("D4").Value = "prep/p03.html"

# Update the active selection to D5 as recorded in the saved workbook
$ws.Range("D5").Select()
